$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "23.0.0"
$ws.Range("A29").Value = "24.0.0"
$ws.Range("A30").Value = "25.0.0"
$ws.Range("A31").Value = "26.0.0"
$ws.Range("A32").Value = "27.0.0"
$ws.Range("B31").Value = "How would you test for this? Maybe this didn’t fulfill SMART "
$ws.Range("B31").IndentLevel = 0
$ws.Range("A33").Value = "28.0.0"
$ws.Range("A34").Value = "29.0.0"
$ws.Range("A35").Value = "30.0.0"
$ws.Range("A36").Value = "31.0.0"
$ws.Range("A37").Value = "32.0.0"
$ws.Range("A38").Value = "33.0.0"
$ws.Range("A39").Value = "34.0.0"
$ws.Range("A40").Value = "35.0.0"
$ws.Range("A41").Value = "36.0.0"
$ws.Range("A42").Value = "37.0.0"
$ws.Range("A43").Value = "38.0.0"
$ws.Range("A44").Value = "39.0.0"
$ws.Range("B28").Value = "Full"
$ws.Range("E28").Value = "Full"
$ws.Range("F28").Value = "Full"
$ws.Range("G28").Value = "Full"
$ws.Range("H28").Value = "Full"
$ws.Range("B29").Value = "Full"
$ws.Range("E29").Value = "Full"
$ws.Range("F29").Value = "Full"
$ws.Range("G29").Value = "Full"
$ws.Range("H29").Value = "Full"
$ws.Range("B30").Value = "Full"
$ws.Range("E30").Value = "Full"
$ws.Range("F30").Value = "Full"
$ws.Range("G30").Value = "Full"
$ws.Range("H30").Value = "Full"
$ws.Range("C31").IndentLevel = 0
$ws.Range("D31").IndentLevel = 0
$ws.Range("E31").IndentLevel = 0
$ws.Range("F31").IndentLevel = 0
$ws.Range("G31").IndentLevel = 0
$ws.Range("H31").IndentLevel = 0
$ws.Range("I31").IndentLevel = 0
$ws.Range("J31").IndentLevel = 0
$ws.Range("K31").IndentLevel = 0
$ws.Range("L31").IndentLevel = 0
$ws.Range("M31").IndentLevel = 0
$ws.Range("B32").Value = "Partial"
$ws.Range("B32").IndentLevel = 0
$ws.Range("E32").Value = "Partial"
$ws.Range("E32").IndentLevel = 0
$ws.Range("F32").Value = "Full"
$ws.Range("F32").IndentLevel = 0
$ws.Range("G32").Value = "Full"
$ws.Range("G32").IndentLevel = 0
$ws.Range("H32").Value = "Full"
$ws.Range("H32").IndentLevel = 0
$ws.Range("B33").Value = "Full"
$ws.Range("B33").IndentLevel = 0
$ws.Range("C33").Value = "Full"
$ws.Range("F33").Value = "Full"
$ws.Range("F33").IndentLevel = 0
$ws.Range("G33").Value = "Full"
$ws.Range("G33").IndentLevel = 0
$ws.Range("H33").Value = "Full"
$ws.Range("H33").IndentLevel = 0
$ws.Range("B34").Value = "Full"
$ws.Range("B34").IndentLevel = 0
$ws.Range("C34").Value = "Full"
$ws.Range("F34").Value = "Full"
$ws.Range("F34").IndentLevel = 0
$ws.Range("G34").Value = "Full"
$ws.Range("G34").IndentLevel = 0
$ws.Range("H34").Value = "Full"
$ws.Range("H34").IndentLevel = 0
$ws.Range("B35").Value = "Full"
$ws.Range("B35").IndentLevel = 0
$ws.Range("C35").Value = "Full"
$ws.Range("F35").Value = "Full"
$ws.Range("F35").IndentLevel = 0
$ws.Range("G35").Value = "Full"
$ws.Range("G35").IndentLevel = 0
$ws.Range("H35").Value = "Full"
$ws.Range("H35").IndentLevel = 0
$ws.Range("B36").Value = "How would you test for this? Maybe this didn’t fulfill SMART "
$ws.Range("B36").IndentLevel = 0
$ws.Range("B37").Value = "Partial "
$ws.Range("B37").IndentLevel = 0
$ws.Range("C37").Value = "Partial"
$ws.Range("F37").Value = "Partial"
$ws.Range("G37").Value = "Partial"
$ws.Range("H37").Value = "Partial"
$ws.Range("B38").Value = "How would you test for this? Maybe this didn’t fulfill SMART "
$ws.Range("B38").IndentLevel = 0
$ws.Range("B39").Value = "Full"
$ws.Range("B39").IndentLevel = 0
$ws.Range("D39").Value = "Full"
$ws.Range("F39").Value = "Full"
$ws.Range("G39").Value = "Full"
$ws.Range("H39").Value = "Full"
$ws.Range("F40").Value = "Full"
$ws.Range("G40").Value = "Full"
$ws.Range("B42").Value = "Full"
$ws.Range("C42").Value = "Full"
$ws.Range("E42").Value = "Full"
$ws.Range("F42").Value = "Full"
$ws.Range("G42").Value = "Full"
$ws.Range("H42").Value = "Full"
$ws.Range("B43").Value = "Full"
$ws.Range("C43").Value = "Full"
$ws.Range("D43").Value = "Full"
$ws.Range("E43").Value = "Full"
$ws.Range("F43").Value = "Full"
$ws.Range("G43").Value = "Full"
$ws.Range("H43").Value = "Full"
$ws.Range("B44").Value = "Full"
$ws.Range("D44").Value = "Full"
$ws.Range("F44").Value = "Full"
$ws.Range("G44").Value = "Full"
$ws.Range("H44").Value = "Full"

$ws.Range("K45").Select() | Out-Null
